$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.942.29"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -4.75%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.221.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -6.17%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'316.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +1.77%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'99.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -7.99%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  -6.44%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  -8.24%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'36.99"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -9.32%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'53.99"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -2.81%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.0829"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -9.40%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'7.79"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -7.54%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D15").Value = "'0.862"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -11.50%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'2.556.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -6.43%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'14.24"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -6.35%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.206.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -7.31%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'42.866.73"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -4.91%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'15.46"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +6.91%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.0₃0965"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -8.91%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'6.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -10.89%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'65.43"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -10.64%  "
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'  -8.70%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'236.85"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -8.54%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'2.12"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -7.96%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  -0.05%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'10.10"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -9.38%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'2.22"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -5.13%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'6.36"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -11.37%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'20.51"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -8.08%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'0.0883"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -8.41%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'34.21"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -8.15%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'157.49"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -6.59%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  -6.02%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'3.21"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +9.80%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'1.98"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +13.02%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = "'  -5.82%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'4.46"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -5.46%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'3.78"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -3.37%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  -11.23%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.0326"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -7.65%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'1.915.20"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.99%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  +0.02%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'12.34"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -3.31%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'89.31"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -10.53%  "
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'  -9.26%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'5.42"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -3.49%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'60.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -12.39%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'75.51"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -6.68%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.869"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +18.12%  "
$ws.Range("E51").ClearFormats()
